# Natmi following Dr Hou advice
#
# The NATMI ligand-receptor table (Cdh1 -> Itgb7) is recomputed with a new
# sending/target cluster "sCs" added alongside the existing "M1"/"M2"
# clusters, and the derived statistics for every sending/target cluster
# pair are refreshed to match the new run. Rows 2-11 (M1/M2 as sender)
# get updated numbers, and rows 12-16 (the new "sCs" sender) are appended.
#
# Columns: A=Sending cluster, B=Ligand symbol, C=Receptor symbol,
#          D=Target cluster, E=Ligand-expressing cells, F=Ligand detection
#          rate, G=Ligand average expression value, H=Ligand total
#          expression value, I=Ligand derived specificity (avg),
#          J=Ligand derived specificity (total), K=Receptor-expressing
#          cells, L=Receptor detection rate, M=Receptor average expression
#          value, N=Receptor total expression value, O=Receptor derived
#          specificity (avg), P=Receptor derived specificity (total),
#          Q=Edge average expression weight, R=Edge total expression
#          weight, S=Edge average expression derived specificity,
#          T=Edge total expression derived specificity.
$data = @(
    @("M1", "Cdh1", "Itgb7", "ECs", 3, 1, 0.4197659999999999, 1.259298, 0.4104130959562529, 0.4104130959562529, 3, 1, 1.570106666666667, 4.71032, 0.04895814174088396, 0.04895814174088397, 0.65907739504, 5.931696555359999, 0.02009306252414124, 0.02009306252414124),
    @("M1", "Cdh1", "Itgb7", "FAPs", 3, 1, 0.4197659999999999, 1.259298, 0.4104130959562529, 0.4104130959562529, 3, 1, 2.170845333333333, 6.512536, 0.06769002118340356, 0.06769002118340357, 0.9112470621919998, 8.201223559727998, 0.027780871159225, 0.027780871159225),
    @("M1", "Cdh1", "Itgb7", "M1", 3, 1, 0.4197659999999999, 1.259298, 0.4104130959562529, 0.4104130959562529, 3, 1, 12.02192266666667, 36.06576800000001, 0.3748605151535007, 0.3748605151535007, 5.046394390096, 45.417549510864, 0.1538476645759041, 0.1538476645759041),
    @("M1", "Cdh1", "Itgb7", "M2", 3, 1, 0.4197659999999999, 1.259298, 0.4104130959562529, 0.4104130959562529, 3, 1, 15.862603, 47.587809, 0.4946183482566181, 0.4946183482566182, 6.658581410897999, 59.92723269808199, 0.2029978476247667, 0.2029978476247668),
    @("M1", "Cdh1", "Itgb7", "sCs", 3, 1, 0.4197659999999999, 1.259298, 0.4104130959562529, 0.4104130959562529, 3, 1, 0.4449116666666666, 1.334735, 0.01387297366559358, 0.01387297366559358, 0.18675879067, 1.68082911603, 0.005693650072215828, 0.005693650072215829),
    @("M2", "Cdh1", "Itgb7", "ECs", 3, 1, 0.5921903333333333, 1.776571, 0.5789956025468922, 0.5789956025468922, 3, 1, 1.570106666666667, 4.71032, 0.04895814174088396, 0.04895814174088397, 0.9298019903022222, 8.36821791272, 0.02834654877683926, 0.02834654877683927),
    @("M2", "Cdh1", "Itgb7", "FAPs", 3, 1, 0.5921903333333333, 1.776571, 0.5789956025468922, 0.5789956025468922, 3, 1, 2.170845333333333, 6.512536, 0.06769002118340356, 0.06769002118340357, 1.285553621561778, 11.569982594056, 0.03919222460149664, 0.03919222460149665),
    @("M2", "Cdh1", "Itgb7", "M1", 3, 1, 0.5921903333333333, 1.776571, 0.5789956025468922, 0.5789956025468922, 3, 1, 12.02192266666667, 36.06576800000001, 0.3748605151535007, 0.3748605151535007, 7.119266391280889, 64.07339752152801, 0.2170425898423395, 0.2170425898423396),
    @("M2", "Cdh1", "Itgb7", "M2", 3, 1, 0.5921903333333333, 1.776571, 0.5789956025468922, 0.5789956025468922, 3, 1, 15.862603, 47.587809, 0.4946183482566181, 0.4946183482566182, 9.393680158104333, 84.543121422939, 0.2863818485795891, 0.2863818485795892),
    @("M2", "Cdh1", "Itgb7", "sCs", 3, 1, 0.5921903333333333, 1.776571, 0.5789956025468922, 0.5789956025468922, 3, 1, 0.4449116666666666, 1.334735, 0.01387297366559358, 0.01387297366559358, 0.2634723881872222, 2.371251493685, 0.008032390746627522, 0.008032390746627524),
    @("sCs", "Cdh1", "Itgb7", "ECs", 1, 0.3333333333333333, 0.01083266666666667, 0.032498, 0.01059130149685484, 0.01059130149685484, 3, 1, 1.570106666666667, 4.71032, 0.04895814174088396, 0.04895814174088397, 0.01700844215111111, 0.15307597936, 0.0005185304399034558, 0.0005185304399034559),
    @("sCs", "Cdh1", "Itgb7", "FAPs", 1, 0.3333333333333333, 0.01083266666666667, 0.032498, 0.01059130149685484, 0.01059130149685484, 3, 1, 2.170845333333333, 6.512536, 0.06769002118340356, 0.06769002118340357, 0.02351604388088889, 0.211644394928, 0.000716925422681918, 0.0007169254226819182),
    @("sCs", "Cdh1", "Itgb7", "M1", 1, 0.3333333333333333, 0.01083266666666667, 0.032498, 0.01059130149685484, 0.01059130149685484, 3, 1, 12.02192266666667, 36.06576800000001, 0.3748605151535007, 0.3748605151535007, 0.1302294809404445, 1.172065328464, 0.003970260735257049, 0.003970260735257049),
    @("sCs", "Cdh1", "Itgb7", "M2", 1, 0.3333333333333333, 0.01083266666666667, 0.032498, 0.01059130149685484, 0.01059130149685484, 3, 1, 15.862603, 47.587809, 0.4946183482566181, 0.4946183482566182, 0.1718342907646667, 1.546508616882, 0.005238652052262188, 0.005238652052262189),
    @("sCs", "Cdh1", "Itgb7", "sCs", 1, 0.3333333333333333, 0.01083266666666667, 0.032498, 0.01059130149685484, 0.01059130149685484, 3, 1, 0.4449116666666666, 1.334735, 0.01387297366559358, 0.01387297366559358, 0.004819579781111111, 0.04337621803, 0.0001469328467502291, 0.0001469328467502291),
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data starts on row 2 (row 1 is the header). Rows 2-11 already existed
# and are overwritten in place; rows 12-16 are new and extend the sheet's
# used range to A1:T16.
$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $row.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
